$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Espinaca" (Vega Modelo de Temuco) is
# inserted as row 35, pushing the existing rows 35-80 down to 36-81.
$ws.Rows(35).Insert()

$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44467
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 100112012
$ws.Range("G35").Value = "Espinaca"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 40
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 11000
$ws.Range("M35").Value = 10500
$ws.Range("N35").Value = "`$/docena de atados"
$ws.Range("O35").Value = "Región de La Araucanía"
$ws.Range("P35").Value = 3500
$ws.Range("Q35").Value = 3
$ws.Range("R35").Value = "Hortaliza"
